$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8826526999473572
$ws.Range("B1").Value = 2.814315795898438
$ws.Range("C1").Value = 4.651193141937256
$ws.Range("D1").Value = 2.365373849868774
$ws.Range("E1").Value = 1.068786978721619
